# The source data table listed four years of data (2008-2011) in rows 2-5.
# The update removes the two oldest years (2008年, 2009年) from the top of
# the table, so 2010年 and 2011年 shift up into rows 2 and 3 and the sheet
# shrinks from A1:U5 to A1:U3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting row 2 twice removes the original 2008年 row, then the original
# 2009年 row (which slides into row 2 after the first delete), shifting the
# remaining 2010年/2011年 rows up to rows 2 and 3.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
